# Foi ajustado o erro na parte de edição na classe dermo
# Adds a new data row (row 2) to the worksheet, matching the header
# columns Data, Meta, Meta.AC, Venda, Venda.AC, Sobras, P.
# The values are written as literal text (leading apostrophe forces a
# text/quoted entry) so numeric- and date-looking strings such as
# "05/05/5000" or "5000.00" are stored verbatim instead of being
# reinterpreted as a date serial or a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "'05/05/5000"
$ws.Range("B2").Value = "'5000.00"
$ws.Range("C2").Value = "'5000.00"
$ws.Range("D2").Value = "'5000.00"
$ws.Range("E2").Value = "'5000.00"
$ws.Range("F2").Value = "'0.00"
$ws.Range("G2").Value = "'100.00"
